$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3624364103015721
$ws.Range("C2").Value = 2.902560095337955
$ws.Range("D2").Value = 14.15319625081184
$ws.Range("E2").Value = 3.762073397850159
$ws.Range("F2").Value = 3.78183458407233
$ws.Range("G2").Value = 51
$ws.Range("B3").Value = 0.03493740399305523
$ws.Range("C3").Value = 3.008365163207992
$ws.Range("D3").Value = 12.2517222994939
$ws.Range("E3").Value = 3.500246034137301
$ws.Range("F3").Value = 3.535606300991146
$ws.Range("G3").Value = 50
$ws.Range("B4").Value = 0.376526435448383
$ws.Range("C4").Value = 2.512159032405139
$ws.Range("D4").Value = 10.14738472548122
$ws.Range("E4").Value = 3.185495993637603
$ws.Range("F4").Value = 3.195944747761078
$ws.Range("G4").Value = 49
$ws.Range("B5").Value = 0.128062596573612
$ws.Range("C5").Value = 2.051014053753879
$ws.Range("D5").Value = 8.415776543719975
$ws.Range("E5").Value = 2.900995784850432
$ws.Range("F5").Value = 2.92883708213644
$ws.Range("G5").Value = 48
$ws.Range("B6").Value = 0.3895087488094023
$ws.Range("C6").Value = 2.65705602853736
$ws.Range("D6").Value = 11.08899265343494
$ws.Range("E6").Value = 3.330013911898108
$ws.Range("F6").Value = 3.342909279152719
$ws.Range("G6").Value = 47
$ws.Range("B7").Value = 0.1522837063575818
$ws.Range("C7").Value = 2.670556882397438
$ws.Range("D7").Value = 10.36728199909625
$ws.Range("E7").Value = 3.219826392695148
$ws.Range("F7").Value = 3.251762656728451
$ws.Range("G7").Value = 46
$ws.Range("B8").Value = 0.3987337555141685
$ws.Range("C8").Value = 2.763470879128723
$ws.Range("D8").Value = 12.14087511611534
$ws.Range("E8").Value = 3.48437585746936
$ws.Range("F8").Value = 3.500600250070111
$ws.Range("G8").Value = 45
$ws.Range("B9").Value = 0.1063170287658169
$ws.Range("C9").Value = 2.290045569855232
$ws.Range("D9").Value = 9.722226614463244
$ws.Range("E9").Value = 3.118048526637013
$ws.Range("F9").Value = 3.152262477892468
$ws.Range("G9").Value = 44
$ws.Range("B10").Value = 0.458157375663734
$ws.Range("C10").Value = 2.282092575471136
$ws.Range("D10").Value = 8.777379630194348
$ws.Range("E10").Value = 2.962664279022237
$ws.Range("F10").Value = 2.961664880566207
$ws.Range("G10").Value = 43
$ws.Range("B11").Value = 0.06186012078826768
$ws.Range("C11").Value = 2.779572087388781
$ws.Range("D11").Value = 10.64187340981079
$ws.Range("E11").Value = 3.262188438734156
$ws.Range("F11").Value = 3.301137878033202
$ws.Range("G11").Value = 42
